$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 322398.8857866722
$ws.Range("D2").Value = 0.2754491017964072
$ws.Range("E2").Value = 0.2765531062124248
$ws.Range("F2").Value = 0.276
$ws.Range("G2").Value = 0.0008577945278104697
$ws.Range("C3").Value = 51678890.85532195
$ws.Range("D3").Value = 0.4037209302325581
$ws.Range("E3").Value = 0.4049586776859504
$ws.Range("F3").Value = 0.4043388567245624
$ws.Range("G3").Value = 0.00005887606809074921
$ws.Range("C4").Value = 25880917.99959877
$ws.Range("D4").Value = 0.3803761673023807
$ws.Range("E4").Value = 0.3854972007464676
$ws.Range("F4").Value = 0.3829195630585899
$ws.Range("G4").Value = 0.0001124897284547477
$ws.Range("C5").Value = 51652880.93814024
$ws.Range("E5").Value = 0.7931218341775527
$ws.Range("F5").Value = 0.8846268212905145
$ws.Range("G5").Value = 0.000101901955729869
$ws.Range("C6").Value = 40897708.30791894
$ws.Range("D6").Value = 0.4037209302325581
$ws.Range("E6").Value = 0.4049586776859504
$ws.Range("F6").Value = 0.4043388567245624
$ws.Range("G6").Value = 0.00007439658657507687
$ws.Range("C7").Value = 26011411.68620238
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0.795921087709944
$ws.Range("F7").Value = 0.8863653232390707
$ws.Range("G7").Value = 0.0002034679012776482
$ws.Range("C8").Value = 40877130.14475385
$ws.Range("E8").Value = 0.7931218341775527
$ws.Range("F8").Value = 0.8846268212905145
$ws.Range("G8").Value = 0.0001287646556409264
$ws.Range("C9").Value = 20487286.90503867
$ws.Range("D9").Value = 0.3803761673023807
$ws.Range("E9").Value = 0.3854972007464676
$ws.Range("F9").Value = 0.3829195630585899
$ws.Range("G9").Value = 0.0001421045866848499
$ws.Range("C10").Value = 20590528.84795526
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0.795921087709944
$ws.Range("F10").Value = 0.8863653232390707
$ws.Range("G10").Value = 0.0002570350370377233
